$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("tesla", 386000000),
    @("question", 4050000000),
    @("register", 2850000000),
    @("complete", 2880000000),
    @("tesla", 375000000),
    @("mridangam", 1110000),
    @("clarity", 210000000),
    @("shoot", 575000000),
    @("tesla", 332000000),
    @("shoot", 4780000000),
    @("cheer", 196000000),
    @("small", 5590000000)
)

$startRow = 75
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
